# payment_links_excel.xlsx — "Add files via upload"
#
# Adds two new ticket-holder rows (Vaasu Bisht / Vansh Duggar, "Omdena
# Event") to sheet1, corrects two phone numbers (now stored as real
# numbers rather than text) and fixes Rishabh Kabra's registration
# number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Existing-row fixups
# ---------------------------------------------------------------------

# Row 2 (Anand Lahoti): corrected phone number, now a real number.
$ws.Range("N2").Value = 9893353361

# Row 3 (Rishabh Kabra): corrected phone number (now numeric) and a
# corrected registration number.
$ws.Range("N3").Value = 9310370781
$ws.Range("P3").Value = "21MIZ10035"

# ---------------------------------------------------------------------
# New row 4 — Vaasu Bisht (Omdena Event)
#
# Columns that repeat a value already present elsewhere in the sheet
# are filled in with Range.Copy so the existing shared-string entry
# (and default "no extra style") is reused verbatim, instead of
# letting the COM layer reinterpret a numeric-looking literal (e.g.
# "1.00") as a true number.
# ---------------------------------------------------------------------

$ws.Range("A2").Copy($ws.Range("A4"))
$ws.Range("B4").Value = "Omdena Event"
$ws.Range("C2").Copy($ws.Range("C4"))
$ws.Range("D4").Value = "order_LQf8k4C4lj53423"
$ws.Range("E2").Copy($ws.Range("E4"))
$ws.Range("F2").Copy($ws.Range("F4"))
$ws.Range("G2").Copy($ws.Range("G4"))
$ws.Range("H2").Copy($ws.Range("H4"))
$ws.Range("I2").Copy($ws.Range("I4"))
$ws.Range("J2").Copy($ws.Range("J4"))
$ws.Range("K3").Copy($ws.Range("K4"))
$ws.Range("L2").Copy($ws.Range("L4"))
$ws.Hyperlinks.Add($ws.Range("M4"), "mailto:rishabh.kabra2021@vitbhopal.ac.in", "", "", "rishabh.kabra2021@vitbhopal.ac.in")
$ws.Range("M4").Style = "Hyperlink"
$ws.Range("N4").Value = 9893323161
$ws.Range("O4").Value = "Vaasu Bisht"
$ws.Range("P4").Value = "21MIM10035"
$ws.Range("Q2").Copy($ws.Range("Q4"))

# ---------------------------------------------------------------------
# New row 5 — Vansh Duggar (Omdena Event)
# ---------------------------------------------------------------------

$ws.Range("A2").Copy($ws.Range("A5"))
$ws.Range("B4").Copy($ws.Range("B5"))
$ws.Range("C2").Copy($ws.Range("C5"))
$ws.Range("D5").Value = "order_LQf8k4C4lj5b327"
$ws.Range("E2").Copy($ws.Range("E5"))
$ws.Range("F2").Copy($ws.Range("F5"))
$ws.Range("G2").Copy($ws.Range("G5"))
$ws.Range("H2").Copy($ws.Range("H5"))
$ws.Range("I2").Copy($ws.Range("I5"))
$ws.Range("J2").Copy($ws.Range("J5"))
$ws.Range("K3").Copy($ws.Range("K5"))
$ws.Range("L2").Copy($ws.Range("L5"))
$ws.Hyperlinks.Add($ws.Range("M5"), "mailto:rishabh.kabra2021@vitbhopal.ac.in", "", "", "rishabh.kabra2021@vitbhopal.ac.in")
$ws.Range("M5").Style = "Hyperlink"
$ws.Range("N5").Value = 9893768161
$ws.Range("O5").Value = "Vansh Duggar"
$ws.Range("P5").Value = "21MIM10039"
$ws.Range("Q2").Copy($ws.Range("Q5"))

# ---------------------------------------------------------------------
# Cosmetic view updates (zoom + selection), matching the re-saved file
# ---------------------------------------------------------------------

$ws.Range("O11").Select()
$excel.ActiveWindow.Zoom = 175

# Column widths, approximating the re-saved layout.
$ws.Columns.Item(1).ColumnWidth = 29.54
$ws.Columns.Item(2).ColumnWidth = 18.09
$ws.Columns.Item(3).ColumnWidth = 30.63
$ws.Columns.Item(4).ColumnWidth = 37.63
$ws.Columns.Item(13).ColumnWidth = 37.45
$ws.Columns.Item(14).ColumnWidth = 19.18
$ws.Columns.Item(15).ColumnWidth = 14.82
$ws.Columns.Item(16).ColumnWidth = 17.54
